# Update the "Date" column (B2:B29) on Sheet1 with a new batch of execution
# timestamps, matching the author's RAD test-data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDates = @(
    "Mon Oct 09 22:32:29 EDT 2023",
    "Mon Oct 09 22:32:39 EDT 2023",
    "Mon Oct 09 22:32:49 EDT 2023",
    "Mon Oct 09 22:32:59 EDT 2023",
    "Mon Oct 09 22:33:09 EDT 2023",
    "Mon Oct 09 22:33:19 EDT 2023",
    "Mon Oct 09 22:33:30 EDT 2023",
    "Mon Oct 09 22:33:40 EDT 2023",
    "Mon Oct 09 22:33:50 EDT 2023",
    "Mon Oct 09 22:34:01 EDT 2023",
    "Mon Oct 09 22:34:11 EDT 2023",
    "Mon Oct 09 22:34:21 EDT 2023",
    "Mon Oct 09 22:34:31 EDT 2023",
    "Mon Oct 09 22:34:41 EDT 2023",
    "Mon Oct 09 22:34:52 EDT 2023",
    "Mon Oct 09 22:35:02 EDT 2023",
    "Mon Oct 09 22:35:12 EDT 2023",
    "Mon Oct 09 22:35:21 EDT 2023",
    "Mon Oct 09 22:35:31 EDT 2023",
    "Mon Oct 09 22:35:41 EDT 2023",
    "Mon Oct 09 22:35:52 EDT 2023",
    "Mon Oct 09 22:36:01 EDT 2023",
    "Mon Oct 09 22:36:11 EDT 2023",
    "Mon Oct 09 22:36:21 EDT 2023",
    "Mon Oct 09 22:36:31 EDT 2023",
    "Mon Oct 09 22:36:42 EDT 2023",
    "Mon Oct 09 22:36:52 EDT 2023",
    "Mon Oct 09 22:37:02 EDT 2023"
)

$startRow = 2
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
}
